# Controls / Camera work in most cases - Strange behaviour on RESET
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raw mark values that feed the weighted-average formulas
$ws.Range("D8").Value = 7
$ws.Range("D10").Value = 0
$ws.Range("D14").Value = 7
$ws.Range("D15").Value = 7
$ws.Range("D35").Value = 10
$ws.Range("D36").Value = 10

# Update the last active cell / selection on the sheet view
$ws.Range("F26").Select()
